# Weekly fruit/veg price update: insert two new observation rows for the
# Vega Modelo de Temuco - Cilantro subset, pushing the existing rows
# (522:638) down by two positions (to 524:640) and populating the two
# newly inserted rows (522:523) with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 522:523 - this shifts rows 522:638 down to
# 524:640, carrying their values/formats with them automatically.
$ws.Rows("522:523").Insert()

# --- New row 522 ---
$ws.Range("A522").Value = 10
$ws.Range("B522").Value = "Vega Modelo de Temuco"
$ws.Range("C522").Value = "La Araucanía"
$ws.Range("D522").Value2 = 45173
$ws.Range("D522").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E522").Value = 9
$ws.Range("F522").Value = 100112040
$ws.Range("G522").Value = "Cilantro"
$ws.Range("H522").Value = "Sin especificar"
$ws.Range("I522").Value = "Primera"
$ws.Range("J522").Value = 60
$ws.Range("K522").Value = 5000
$ws.Range("L522").Value = 5000
$ws.Range("M522").Value = 5000
$ws.Range("N522").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O522").Value = "Provincia de Cautín"
$ws.Range("P522").Value = 2500
$ws.Range("Q522").Value = 2
$ws.Range("R522").Value = "Hortaliza"

# --- New row 523 ---
$ws.Range("A523").Value = 10
$ws.Range("B523").Value = "Vega Modelo de Temuco"
$ws.Range("C523").Value = "La Araucanía"
$ws.Range("D523").Value2 = 45173
$ws.Range("D523").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E523").Value = 9
$ws.Range("F523").Value = 100112040
$ws.Range("G523").Value = "Cilantro"
$ws.Range("H523").Value = "Sin especificar"
$ws.Range("I523").Value = "Primera"
$ws.Range("J523").Value = 120
$ws.Range("K523").Value = 3300
$ws.Range("L523").Value = 3300
$ws.Range("M523").Value = 3300
$ws.Range("N523").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O523").Value = "Región Metropolitana"
$ws.Range("P523").Value = 1650
$ws.Range("Q523").Value = 2
$ws.Range("R523").Value = "Hortaliza"
